$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D3").Value = -7.895
$ws.Range("C7").Value = -12.844
$ws.Range("A8").Value = -22.188
$ws.Range("A10").Value = -21.897
$ws.Range("E10").Value = 16.516
$ws.Range("A12").Value = -21.52
$ws.Range("E12").Value = 17.499
$ws.Range("E13").Value = 16.416
$ws.Range("E14").Value = 16.913
$ws.Range("C15").Value = -13.728
$ws.Range("A18").Value = -21.862
$ws.Range("C18").Value = -12.294
$ws.Range("D18").Value = -7.889999999999999
$ws.Range("D19").Value = -8.113
$ws.Range("C20").Value = -12.673
$ws.Range("D27").Value = -8.285000000000002
$ws.Range("C29").Value = -12.591
$ws.Range("E29").Value = 16.947
$ws.Range("C30").Value = -12.54
$ws.Range("C31").Value = -12.759
$ws.Range("D31").Value = -8.070000000000002
$ws.Range("E32").Value = 16.797
$ws.Range("E35").Value = 16.358
$ws.Range("A37").Value = -20.029
$ws.Range("D38").Value = -7.726999999999999
$ws.Range("C40").Value = -12.782
$ws.Range("D42").Value = -8.257999999999999
$ws.Range("E43").Value = 16.885
$ws.Range("D44").Value = -7.316000000000001
$ws.Range("D47").Value = -7.394999999999999
$ws.Range("E48").Value = 17.176
$ws.Range("E49").Value = 16.449
$ws.Range("C50").Value = -12.914
$ws.Range("E50").Value = 16.484
$ws.Range("A55").Value = -21.894
$ws.Range("E56").Value = 16.198
$ws.Range("D58").Value = -8.401
$ws.Range("D65").Value = -7.928999999999999
$ws.Range("A68").Value = -21.632
$ws.Range("C68").Value = -10.98
$ws.Range("E69").Value = 17.241
$ws.Range("D73").Value = -7.877000000000001
$ws.Range("C76").Value = -13.045
$ws.Range("A77").Value = -20.519
$ws.Range("A78").Value = -19.951
$ws.Range("A81").Value = -21.795
$ws.Range("E81").Value = 16.67
$ws.Range("A82").Value = -22.152
$ws.Range("C87").Value = -13.193
$ws.Range("C88").Value = -12.877
$ws.Range("D90").Value = -7.452
$ws.Range("E92").Value = 17.637
$ws.Range("D94").Value = -7.462000000000001
$ws.Range("D95").Value = -7.567
$ws.Range("C96").Value = -12.665
$ws.Range("C98").Value = -13.649
$ws.Range("C101").Value = -13.048
$ws.Range("D101").Value = -8.061999999999999
$ws.Range("C102").Value = -13.086
